$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Michael Beaver")

# Use the existing date-formatted cell (A17) as the format source so we
# don't spawn a brand-new numFmt entry in styles.xml.
$ws.Range("A17").Copy()

# Row 18: 2014-02-05, new task, 1 hour
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A18").Value = 41675
$ws.Range("B18").Value = "Created initial draft of specifications document. Began drafting the introductory sections."
$ws.Range("I18").Value = 1

# Row 19: 2014-02-06, new task, 3.5 hours
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("A19").Value = 41676
$ws.Range("B19").Value = "Read IEEE Standard 830-1998 recommendations for specifications documents. Updated specifications document to version 1.0.1 to include definitions, acronyms, abbreviations, update procedures, and appendices."
$ws.Range("I19").Value = 3.5
$ws.Rows.Item(19).RowHeight = 36.75

# Row 20: 2014-02-06, new task, 0.25 hours
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("A20").Value = 41676
$ws.Range("B20").Value = "Impromptu team meeting. Recorded meeting minutes. Set times for next two team meetings."
$ws.Range("I20").Value = 0.25

$excel.CutCopyMode = $false

# Update the sheet view: scroll position and current selection
$ws.Application.ActiveWindow.ScrollRow = 41
$ws.Range("B23:H23").Select()
